$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections (shared-string edits) ---
$ws.Range("B2").Value = "rohan"
$ws.Range("L2").Value = "Reading ,Drawing"
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Row height normalization for the header + two data rows ---
$ws.Rows("1:3").RowHeight = 19.5

# --- Pincode/Phone columns (bordered, number-formatted cells): bake the
#     theme-based text color down to an explicit black, cell by cell so
#     every one of them lands on the same resulting format. ---
foreach ($addr in @("I2", "K2", "I3", "K3")) {
    $ws.Range($addr).Font.Color = 0
}
